$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("J7").Value = -0.0663
$ws.Range("K7").Value = -0.0079
$ws.Range("L7").Value = -0.0126
$ws.Range("M7").Value = -0.0183
$ws.Range("N7").Value = -0.0516
$ws.Range("O7").Value = 0.0033
$ws.Range("P7").Value = 0.0125
$ws.Range("Q7").Value = 0.012
$ws.Range("R7").Value = 0.0213

# Row 9
$ws.Range("H9").Value = -0.0637
$ws.Range("I9").Value = -0.0648
$ws.Range("J9").Value = -0.0475
$ws.Range("K9").Value = -0.1033
$ws.Range("L9").Value = -0.0365
$ws.Range("M9").Value = -0.0239
$ws.Range("N9").Value = -0.0264
$ws.Range("O9").Value = -0.0206
$ws.Range("P9").Value = -0.0218
$ws.Range("Q9").Value = -0.0304
$ws.Range("R9").Value = -0.0213

# Row 16
$ws.Range("H16").Value = -2.486
$ws.Range("I16").Value = -1.0237
$ws.Range("J16").Value = -1.4718
$ws.Range("K16").Value = -1.7899
$ws.Range("L16").Value = -0.8836
$ws.Range("M16").Value = -0.7029
$ws.Range("N16").Value = -0.7425
$ws.Range("O16").Value = -0.3441
$ws.Range("P16").Value = -0.3357
$ws.Range("Q16").Value = -0.0937
$ws.Range("R16").Value = 0.0342

# Row 20
$ws.Range("K20").Value = -0.4285
$ws.Range("L20").Value = -0.0125
$ws.Range("M20").Value = -0.0123
$ws.Range("N20").Value = -0.6088
$ws.Range("O20").Value = -0.0001
$ws.Range("P20").Value = -0.0001
$ws.Range("Q20").Value = -0.0062

# Row 24
$ws.Range("J24").Value = -0.0025
$ws.Range("K24").Value = 0.004
$ws.Range("L24").Value = 0.12
$ws.Range("M24").Value = 0.1118
$ws.Range("N24").Value = 0.0948
$ws.Range("O24").Value = -0.01
$ws.Range("P24").Value = -0.0672
$ws.Range("Q24").Value = -0.055
$ws.Range("R24").Value = -0.0536

# Row 35
$ws.Range("J35").Value = 0.0176
$ws.Range("K35").Value = -0.0715
$ws.Range("L35").Value = 0.0232
$ws.Range("M35").Value = 0.003
$ws.Range("O35").Value = 0.0739
$ws.Range("P35").Value = -0.0224
$ws.Range("Q35").Value = -0.0214
$ws.Range("R35").Value = -0.0205

# Row 37
$ws.Range("H37").Value = 0.0765
$ws.Range("I37").Value = 0.0259
$ws.Range("J37").Value = 0.0149
$ws.Range("K37").Value = -0.0646
$ws.Range("L37").Value = -0.0207
$ws.Range("M37").Value = -0.0061
$ws.Range("N37").Value = 0.013
$ws.Range("O37").Value = 0.0137
$ws.Range("P37").Value = 0.0066
$ws.Range("Q37").Value = -0.0115
$ws.Range("R37").Value = -0.0178

# Row 44
$ws.Range("H44").Value = 0.0769
$ws.Range("I44").Value = 1.078
$ws.Range("J44").Value = 0.3607
$ws.Range("K44").Value = 0.7619
$ws.Range("L44").Value = 0.1807
$ws.Range("M44").Value = -0.2328
$ws.Range("N44").Value = -0.5743
$ws.Range("O44").Value = 0.024
$ws.Range("P44").Value = 0.0026
$ws.Range("Q44").Value = -0.0864
$ws.Range("R44").Value = -0.0811

# Row 48
$ws.Range("K48").Value = 0.6121
$ws.Range("L48").Value = -0.0122
$ws.Range("M48").Value = -0.0121
$ws.Range("N48").Value = -0.5982
$ws.Range("O48").Value = -0.0001
$ws.Range("P48").Value = -0.0001
$ws.Range("Q48").Value = -0.0062

# Row 52
$ws.Range("J52").Value = -0.0074
$ws.Range("K52").Value = 0.0856
$ws.Range("L52").Value = -0.0053
$ws.Range("M52").Value = 0.0009
$ws.Range("N52").Value = -0.002
$ws.Range("O52").Value = -0.0929
$ws.Range("P52").Value = -0.0074
$ws.Range("Q52").Value = -0.0072
$ws.Range("R52").Value = -0.0069
